$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = "'63.017.21"
$ws.Range('E2').Value2 = '  +0.82%  '
$ws.Range('D3').Value2 = "'3.062.14"
$ws.Range('E3').Value2 = '  +0.15%  '
$ws.Range('E4').Value2 = '  -0.03%  '
$ws.Range('D5').Value2 = "'535.63"
$ws.Range('E5').Value2 = '  -0.41%  '
$ws.Range('D6').Value2 = "'136.37"
$ws.Range('E6').Value2 = '  +2.74%  '
$ws.Range('E7').Value2 = '  +0.00%  '
$ws.Range('D8').Value2 = "'3.056.89"
$ws.Range('E8').Value2 = '  +0.20%  '
$ws.Range('E9').Value2 = '  +0.63%  '
$ws.Range('D10').Value2 = "'0.154"
$ws.Range('E10').Value2 = '  +0.87%  '
$ws.Range('D11').Value2 = "'6.18"
$ws.Range('E11').Value2 = '  +0.19%  '
$ws.Range('D12').Value2 = "'0.450"
$ws.Range('E12').Value2 = '  -1.88%  '
$ws.Range('E13').Value2 = '  +0.52%  '
$ws.Range('D14').Value2 = "'34.17"
$ws.Range('E14').Value2 = '  -0.97%  '
$ws.Range('D15').Value2 = "'3.556.42"
$ws.Range('E15').Value2 = '  +0.95%  '
$ws.Range('D16').Value2 = "'62.983.18"
$ws.Range('E16').Value2 = '  +0.76%  '
$ws.Range('D17').Value2 = "'0.113"
$ws.Range('E17').Value2 = '  +1.69%  '
$ws.Range('D18').Value2 = "'3.064.22"
$ws.Range('E18').Value2 = '  +0.09%  '
$ws.Range('D19').Value2 = "'6.58"
$ws.Range('E19').Value2 = '  -0.42%  '
$ws.Range('D20').Value2 = "'467.68"
$ws.Range('E20').Value2 = '  -1.99%  '
$ws.Range('D21').Value2 = "'13.29"
$ws.Range('E21').Value2 = '  -0.29%  '
$ws.Range('D22').Value2 = "'0.689"
$ws.Range('E22').Value2 = '  -1.93%  '
$ws.Range('D23').Value2 = "'6.95"
$ws.Range('E23').Value2 = '  -3.36%  '
$ws.Range('D24').Value2 = "'78.14"
$ws.Range('E24').Value2 = '  -0.20%  '
$ws.Range('D25').Value2 = "'12.02"
$ws.Range('E25').Value2 = '  +0.05%  '
$ws.Range('E26').Value2 = '  +0.09%  '
$ws.Range('D27').Value2 = "'2.67"
$ws.Range('E27').Value2 = '  -1.06%  '
$ws.Range('D28').Value2 = "'7.80"
$ws.Range('E28').Value2 = '  -4.35%  '
$ws.Range('D29').Value2 = "'0.999"
$ws.Range('E29').Value2 = '  -0.09%  '
$ws.Range('D30').Value2 = "'25.93"
$ws.Range('E30').Value2 = '  -0.02%  '
$ws.Range('E31').Value2 = '  +4.50%  '
$ws.Range('E32').Value2 = '  -2.82%  '
$ws.Range('D33').Value2 = "'58.99"
$ws.Range('E33').Value2 = '  +1.17%  '
$ws.Range('E34').Value2 = '  -5.02%  '
$ws.Range('D35').Value2 = "'5.40"
$ws.Range('E35').Value2 = '  +4.96%  '
$ws.Range('E36').Value2 = '  -1.21%  '
$ws.Range('D37').Value2 = "'474.54"
$ws.Range('E37').Value2 = '  -1.79%  '
$ws.Range('D38').Value2 = "'3.229.37"
$ws.Range('E38').Value2 = '  +3.24%  '
$ws.Range('E39').Value2 = '  +1.05%  '
$ws.Range('D40').Value2 = "'0.0786"
$ws.Range('E40').Value2 = '  -0.77%  '
$ws.Range('E41').Value2 = '  +1.01%  '
$ws.Range('D42').Value2 = "'8.06"
$ws.Range('E42').Value2 = '  +0.48%  '
$ws.Range('D43').Value2 = "'2.56"
$ws.Range('E43').Value2 = '  +0.67%  '
$ws.Range('E45').Value2 = '  -1.26%  '
$ws.Range('D46').Value2 = "'123.07"
$ws.Range('E46').Value2 = '  +4.15%  '
$ws.Range('D47').Value2 = "'24.91"
$ws.Range('E47').Value2 = '  +1.91%  '
$ws.Range('E48').Value2 = '  -2.62%  '
$ws.Range('E49').Value2 = '  +1.05%  '
$ws.Range('D50').Value2 = '0.0₃0514'
$ws.Range('E50').Value2 = '  +2.32%  '
$ws.Range('D51').Value2 = "'1.24"
$ws.Range('E51').Value2 = '  +5.35%  '
